$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The phone-number column holds text values (e.g. with leading zeros in the
# source data), so force it to Text format before writing new values in,
# preventing Excel from auto-converting them to numbers.
$ws.Range("B1:B5").NumberFormat = "@"

# Update header row
$ws.Range("B1").Value = "Phone Number"

# Update data rows 2-5 with new contact info
$ws.Range("A2").Value = "john"
$ws.Range("B2").Value = "36478291"
$ws.Range("C2").Value = "john@email.com"

$ws.Range("A3").Value = "brian"
$ws.Range("B3").Value = "2834501"
$ws.Range("C3").Value = "brian@email.com"

$ws.Range("A4").Value = "judith"
$ws.Range("B4").Value = "30981234"
$ws.Range("C4").Value = "judith@email.com"

$ws.Range("A5").Value = "nas"
$ws.Range("B5").Value = "378192304"
$ws.Range("C5").Value = "nas@email.com"

# Remove the old 6th row entirely (shift rows up)
$ws.Rows("6:6").Delete()
